$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 449; this shifts the existing rows 449-512
# down to 450-513 and the sheet dimension grows from R512 to R513
# automatically.
$ws.Rows.Item(449).Insert()

# Populate the newly inserted row 449 with a new data point. The
# non-varying "template" columns (A,B,C,E,F,G,H,I,N,O,Q,R) mirror the
# values already used throughout this block of rows (same market /
# product / region), while D,J,K,L,M,P carry the new record's data.
$ws.Range("A449").Value = 10
$ws.Range("B449").Value = "Vega Modelo de Temuco"
$ws.Range("C449").Value = "La Araucanía"
$ws.Range("D449").Value = 45077
$ws.Range("E449").Value = 9
$ws.Range("F449").Value = 100112009
$ws.Range("G449").Value = "Acelga"
$ws.Range("H449").Value = "Sin especificar"
$ws.Range("I449").Value = "Primera"
$ws.Range("J449").Value = 65
$ws.Range("K449").Value = 8000
$ws.Range("L449").Value = 8000
$ws.Range("M449").Value = 8000
$ws.Range("N449").Value = "$/docena de atados (12 kilos)"
$ws.Range("O449").Value = "Provincia de Cautín"
$ws.Range("P449").Value = 667
$ws.Range("Q449").Value = 12
$ws.Range("R449").Value = "Hortaliza"

# Match the date-number formatting style used by the rest of column D.
$ws.Range("D449").NumberFormat = $ws.Range("D450").NumberFormat
